$wb = $excel.ActiveWorkbook

# --- 1. Update sheet 'Existing username': fix B4 timestamp + append rows 5-8 ---
$wsExisting = $wb.Worksheets.Item("Existing username")
$wsExisting.Range("B4").Value = 45717.90519327547

# row 5
$wsExisting.Range('A5').Value = 'New customer registration'
$wsExisting.Range('B5').Value = 45719.83218118056
$wsExisting.Range('C5').Value = 'Igor'
$wsExisting.Range('D5').Value = 'Moryc'
$wsExisting.Range('E5').Value = 'ulica Kasztanowa 35/17'
$wsExisting.Range('F5').Value = 'Ostrów Mazowiecka'
$wsExisting.Range('G5').Value = 'Lubelskie'
$wsExisting.Range('H5').Value = '39-194'
$wsExisting.Range('I5').Value = '519 130 953'
$wsExisting.Range('J5').NumberFormat = "@"
$wsExisting.Range('J5').Value = '99072765780'
$wsExisting.Range('K5').Value = 'igor_mor'
$wsExisting.Range('L5').Value = '$+40KsFM+$'
$wsExisting.Range("B5").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 6
$wsExisting.Range('A6').Value = 'New customer registration'
$wsExisting.Range('B6').Value = 45719.84338392361
$wsExisting.Range('C6').Value = 'Maksymilian'
$wsExisting.Range('D6').Value = 'Żywica'
$wsExisting.Range('E6').Value = 'ulica Krasickiego 534'
$wsExisting.Range('F6').Value = 'Piekary Śląskie'
$wsExisting.Range('G6').Value = 'Podlaskie'
$wsExisting.Range('H6').Value = '65-632'
$wsExisting.Range('I6').Value = '730 582 732'
$wsExisting.Range('J6').NumberFormat = "@"
$wsExisting.Range('J6').Value = '16252768012'
$wsExisting.Range('K6').Value = 'maksymilian'
$wsExisting.Range('L6').Value = '(V7AHwm%se'
$wsExisting.Range("B6").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 7
$wsExisting.Range('A7').Value = 'New customer registration'
$wsExisting.Range('B7').Value = 45719.84787918982
$wsExisting.Range('C7').Value = 'Jakub'
$wsExisting.Range('D7').Value = 'Pezda'
$wsExisting.Range('E7').Value = 'plac Mazowiecka 88'
$wsExisting.Range('F7').Value = 'Kluczbork'
$wsExisting.Range('G7').Value = 'Podkarpackie'
$wsExisting.Range('H7').Value = '35-903'
$wsExisting.Range('I7').Value = '+48 607 956 141'
$wsExisting.Range('J7').NumberFormat = "@"
$wsExisting.Range('J7').Value = '95012292275'
$wsExisting.Range('K7').Value = 'jakub_pe'
$wsExisting.Range('L7').Value = '(V15YXkv7E'
$wsExisting.Range("B7").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 8
$wsExisting.Range('A8').Value = 'New customer registration'
$wsExisting.Range('B8').Value = 45719.85717928241
$wsExisting.Range('C8').Value = 'Jakub'
$wsExisting.Range('D8').Value = 'Pezda'
$wsExisting.Range('E8').Value = 'plac Mazowiecka 88'
$wsExisting.Range('F8').Value = 'Kluczbork'
$wsExisting.Range('G8').Value = 'Podkarpackie'
$wsExisting.Range('H8').Value = '35-903'
$wsExisting.Range('I8').Value = '+48 607 956 141'
$wsExisting.Range('J8').NumberFormat = "@"
$wsExisting.Range('J8').Value = '95012292275'
$wsExisting.Range('K8').Value = 'jakub_pe'
$wsExisting.Range('L8').Value = '(V15YXkv7E'
$wsExisting.Range("B8").NumberFormat = "yyyy-mm-dd h:mm:ss"

# --- 2. Add new sheet 'Different passwords' (copy layout of 'Existing username') ---
$lastIdx = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIdx)
$wsDiff = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsDiff.Name = "Different passwords"

$wsDiff.Columns.Item(1).ColumnWidth = 30
$wsDiff.Columns.Item(2).ColumnWidth = 30
$wsDiff.Columns.Item(3).ColumnWidth = 30
$wsDiff.Columns.Item(4).ColumnWidth = 30
$wsDiff.Columns.Item(5).ColumnWidth = 30
$wsDiff.Columns.Item(6).ColumnWidth = 30
$wsDiff.Columns.Item(7).ColumnWidth = 30
$wsDiff.Columns.Item(8).ColumnWidth = 30
$wsDiff.Columns.Item(9).ColumnWidth = 30
$wsDiff.Columns.Item(10).ColumnWidth = 30
$wsDiff.Columns.Item(11).ColumnWidth = 30
$wsDiff.Columns.Item(12).ColumnWidth = 30

# header row
$wsDiff.Range('A1').Value = 'New customer registration'
$wsDiff.Range('B1').Value = 'Date & time'
$wsDiff.Range('C1').Value = 'firstname'
$wsDiff.Range('D1').Value = 'lastname'
$wsDiff.Range('E1').Value = 'streetaddress'
$wsDiff.Range('F1').Value = 'city'
$wsDiff.Range('G1').Value = 'state'
$wsDiff.Range('H1').Value = 'postcode'
$wsDiff.Range('I1').Value = 'phonenumber'
$wsDiff.Range('J1').Value = 'ssn'
$wsDiff.Range('K1').Value = 'username'
$wsDiff.Range('L1').Value = 'password'

# row 2
$wsDiff.Range('A2').Value = 'New customer registration'
$wsDiff.Range('B2').Value = 45719.84791469907
$wsDiff.Range('C2').Value = 'Liwia'
$wsDiff.Range('D2').Value = 'Gzyl'
$wsDiff.Range('E2').Value = 'al. Mazowiecka 61'
$wsDiff.Range('F2').Value = 'Jastrzębie-Zdrój'
$wsDiff.Range('G2').Value = 'Pomorskie'
$wsDiff.Range('H2').Value = '64-798'
$wsDiff.Range('I2').Value = '+48 791 593 589'
$wsDiff.Range('J2').NumberFormat = "@"
$wsDiff.Range('J2').Value = '19322164246'
$wsDiff.Range('K2').Value = 'kzywica'
$wsDiff.Range('L2').Value = '4GtN*hv0*Y'
$wsDiff.Range("B2").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 3
$wsDiff.Range('A3').Value = 'New customer registration'
$wsDiff.Range('B3').Value = 45719.85720462963
$wsDiff.Range('C3').Value = 'Maksymilian'
$wsDiff.Range('D3').Value = 'Komisarczyk'
$wsDiff.Range('E3').Value = 'pl. Kolonia 39/34'
$wsDiff.Range('F3').Value = 'Chojnice'
$wsDiff.Range('G3').Value = 'Opolskie'
$wsDiff.Range('H3').Value = '70-122'
$wsDiff.Range('I3').Value = '+48 503 361 238'
$wsDiff.Range('J3').NumberFormat = "@"
$wsDiff.Range('J3').Value = '75102465543'
$wsDiff.Range('K3').Value = 'marcinbosek'
$wsDiff.Range('L3').Value = 'kkRuqNtR#9'
$wsDiff.Range("B3").NumberFormat = "yyyy-mm-dd h:mm:ss"

# --- 3. Add new sheet 'Database verification' ---
$lastIdx2 = $wb.Worksheets.Count
$lastSheet2 = $wb.Worksheets.Item($lastIdx2)
$wsDb = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$wsDb.Name = "Database verification"

$wsDb.Columns.Item(1).ColumnWidth = 30
$wsDb.Columns.Item(2).ColumnWidth = 30
$wsDb.Columns.Item(3).ColumnWidth = 30
$wsDb.Columns.Item(4).ColumnWidth = 30
$wsDb.Columns.Item(5).ColumnWidth = 30
$wsDb.Columns.Item(6).ColumnWidth = 30
$wsDb.Columns.Item(7).ColumnWidth = 30
$wsDb.Columns.Item(8).ColumnWidth = 30
$wsDb.Columns.Item(9).ColumnWidth = 30
$wsDb.Columns.Item(10).ColumnWidth = 30
$wsDb.Columns.Item(11).ColumnWidth = 30
$wsDb.Columns.Item(12).ColumnWidth = 30

# header row
$wsDb.Range('A1').Value = 'New customer registration'
$wsDb.Range('B1').Value = 'Date & time'
$wsDb.Range('C1').Value = 'firstname'
$wsDb.Range('D1').Value = 'lastname'
$wsDb.Range('E1').Value = 'streetaddress'
$wsDb.Range('F1').Value = 'city'
$wsDb.Range('G1').Value = 'state'
$wsDb.Range('H1').Value = 'postcode'
$wsDb.Range('I1').Value = 'phonenumber'
$wsDb.Range('J1').Value = 'ssn'
$wsDb.Range('K1').Value = 'username'
$wsDb.Range('L1').Value = 'password'

# row 2
$wsDb.Range('A2').Value = 'New customer registration'
$wsDb.Range('B2').Value = 45719.85760917824
$wsDb.Range('C2').Value = 'Tymon'
$wsDb.Range('D2').Value = 'Sabała'
$wsDb.Range('E2').Value = 'ulica Bydgoska 550'
$wsDb.Range('F2').Value = 'Kielce'
$wsDb.Range('G2').Value = 'Zachodniopomorskie'
$wsDb.Range('H2').Value = '02-009'
$wsDb.Range('I2').Value = '+48 733 510 413'
$wsDb.Range('J2').NumberFormat = "@"
$wsDb.Range('J2').Value = '12292369578'
$wsDb.Range('K2').Value = 'arkadiuszrzezniczek'
$wsDb.Range('L2').Value = 'QWF6fsjk&3'
$wsDb.Range("B2").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 3
$wsDb.Range('A3').Value = 'New customer registration'
$wsDb.Range('B3').Value = 45719.85809120371
$wsDb.Range('C3').Value = 'Maciej'
$wsDb.Range('D3').Value = 'Kurcz'
$wsDb.Range('E3').Value = 'al. Broniewskiego 63'
$wsDb.Range('F3').Value = 'Kluczbork'
$wsDb.Range('G3').Value = 'Wielkopolskie'
$wsDb.Range('H3').Value = '32-906'
$wsDb.Range('I3').Value = '507 015 051'
$wsDb.Range('J3').NumberFormat = "@"
$wsDb.Range('J3').Value = '78112951369'
$wsDb.Range('K3').Value = 'oskarjonczyk'
$wsDb.Range('L3').Value = 'AV!g5FAxc4'
$wsDb.Range("B3").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 4
$wsDb.Range('A4').Value = 'New customer registration'
$wsDb.Range('B4').Value = 45719.86011217593
$wsDb.Range('C4').Value = 'Radosław'
$wsDb.Range('D4').Value = 'Wydmuch'
$wsDb.Range('E4').Value = 'aleja Dworska 426'
$wsDb.Range('F4').Value = 'Swarzędz'
$wsDb.Range('G4').Value = 'Mazowieckie'
$wsDb.Range('H4').Value = '39-103'
$wsDb.Range('I4').Value = '+48 789 740 913'
$wsDb.Range('J4').NumberFormat = "@"
$wsDb.Range('J4').Value = '24210649228'
$wsDb.Range('K4').Value = 'marikaklaja'
$wsDb.Range('L4').Value = '%^3YXwHj4D'
$wsDb.Range("B4").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 5
$wsDb.Range('A5').Value = 'New customer registration'
$wsDb.Range('B5').Value = 45719.86066746528
$wsDb.Range('C5').Value = 'Bianka'
$wsDb.Range('D5').Value = 'Sapała'
$wsDb.Range('E5').Value = 'aleja Waryńskiego 548'
$wsDb.Range('F5').Value = 'Świętochłowice'
$wsDb.Range('G5').Value = 'Śląskie'
$wsDb.Range('H5').Value = '58-775'
$wsDb.Range('I5').Value = '+48 22 508 09 00'
$wsDb.Range('J5').NumberFormat = "@"
$wsDb.Range('J5').Value = '81112408915'
$wsDb.Range('K5').Value = 'upeksa'
$wsDb.Range('L5').Value = 'b9#MbuFX(V'
$wsDb.Range("B5").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 6
$wsDb.Range('A6').Value = 'New customer registration'
$wsDb.Range('B6').Value = 45719.86157180556
$wsDb.Range('C6').Value = 'Jędrzej'
$wsDb.Range('D6').Value = 'Mrózek'
$wsDb.Range('E6').Value = 'ul. Reymonta 19/56'
$wsDb.Range('F6').Value = 'Pruszcz Gdański'
$wsDb.Range('G6').Value = 'Pomorskie'
$wsDb.Range('H6').Value = '23-434'
$wsDb.Range('I6').Value = '+48 32 039 64 36'
$wsDb.Range('J6').NumberFormat = "@"
$wsDb.Range('J6').Value = '99071546751'
$wsDb.Range('K6').Value = 'nkuban'
$wsDb.Range('L6').Value = ')2DWkMhFV8'
$wsDb.Range("B6").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 7
$wsDb.Range('A7').Value = 'New customer registration'
$wsDb.Range('B7').Value = 45719.86267232639
$wsDb.Range('C7').Value = 'Kamila'
$wsDb.Range('D7').Value = 'Wawrzynowicz'
$wsDb.Range('E7').Value = 'pl. Południowa 212'
$wsDb.Range('F7').Value = 'Świętochłowice'
$wsDb.Range('G7').Value = 'Lubelskie'
$wsDb.Range('H7').Value = '05-777'
$wsDb.Range('I7').Value = '794 362 300'
$wsDb.Range('J7').NumberFormat = "@"
$wsDb.Range('J7').Value = '86122582770'
$wsDb.Range('K7').Value = 'karina06'
$wsDb.Range('L7').Value = '@lISCIkwx3'
$wsDb.Range("B7").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 8
$wsDb.Range('A8').Value = 'New customer registration'
$wsDb.Range('B8').Value = 45719.86440782408
$wsDb.Range('C8').Value = 'Borys'
$wsDb.Range('D8').Value = 'Siara'
$wsDb.Range('E8').Value = 'ulica Długa 09'
$wsDb.Range('F8').Value = 'Krotoszyn'
$wsDb.Range('G8').Value = 'Kujawsko - pomorskie'
$wsDb.Range('H8').Value = '41-880'
$wsDb.Range('I8').Value = '+48 720 718 425'
$wsDb.Range('J8').NumberFormat = "@"
$wsDb.Range('J8').Value = '06242147242'
$wsDb.Range('K8').Value = 'ada14'
$wsDb.Range('L8').Value = '+2TCa3H0P3'
$wsDb.Range("B8").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 9
$wsDb.Range('A9').Value = 'New customer registration'
$wsDb.Range('B9').Value = 45719.86669033565
$wsDb.Range('C9').Value = 'Tomasz'
$wsDb.Range('D9').Value = 'Cieciura'
$wsDb.Range('E9').Value = 'al. Bursztynowa 73'
$wsDb.Range('F9').Value = 'Chorzów'
$wsDb.Range('G9').Value = 'Kujawsko - pomorskie'
$wsDb.Range('H9').Value = '38-993'
$wsDb.Range('I9').Value = '734 451 473'
$wsDb.Range('J9').NumberFormat = "@"
$wsDb.Range('J9').Value = '00251209133'
$wsDb.Range('K9').Value = 'sebastian82'
$wsDb.Range('L9').Value = 'w_2L+aIr4A'
$wsDb.Range("B9").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 10
$wsDb.Range('A10').Value = 'New customer registration'
$wsDb.Range('B10').Value = 45719.89587921296
$wsDb.Range('C10').Value = 'Daniel'
$wsDb.Range('D10').Value = 'Potęga'
$wsDb.Range('E10').Value = 'pl. Brzechwy 87'
$wsDb.Range('F10').Value = 'Piekary Śląskie'
$wsDb.Range('G10').Value = 'Dolnośląskie'
$wsDb.Range('H10').Value = '65-377'
$wsDb.Range('I10').Value = '603 449 853'
$wsDb.Range('J10').NumberFormat = "@"
$wsDb.Range('J10').Value = '17271641395'
$wsDb.Range('K10').Value = 'ingaczura'
$wsDb.Range('L10').Value = 'E8^2+CNmnt'
$wsDb.Range("B10").NumberFormat = "yyyy-mm-dd h:mm:ss"

# row 11
$wsDb.Range('A11').Value = 'New customer registration'
$wsDb.Range('B11').Value = 45719.89811395122
$wsDb.Range('C11').Value = 'Sonia'
$wsDb.Range('D11').Value = 'Sporysz'
$wsDb.Range('E11').Value = 'aleja Odrzańska 81/28'
$wsDb.Range('F11').Value = 'Gniezno'
$wsDb.Range('G11').Value = 'Lubuskie'
$wsDb.Range('H11').Value = '75-530'
$wsDb.Range('I11').Value = '607 196 066'
$wsDb.Range('J11').NumberFormat = "@"
$wsDb.Range('J11').Value = '98040288113'
$wsDb.Range('K11').Value = 'bielasstanislaw'
$wsDb.Range('L11').Value = '((a2UYh_f^'
$wsDb.Range("B11").NumberFormat = "yyyy-mm-dd h:mm:ss"

Write-Output "done"
